# EPEX spot / Gaz / CO2 daily data refresh
# ------------------------------------------------------------------
# 1) "Prix Spot" sheet: a new daily price column for "01-nov" is
#    inserted right before the "01-oct." column (currently column DF),
#    pushing every existing October/November column one slot to the
#    right (DF->DG, DG->DH, ... EJ->EK). The freshly inserted column
#    gets the "01-nov" header and "-" placeholders for every hour row
#    (no data yet for that day).
# 2) "Gaz" and "CO2" sheets each get one new trailing row with the
#    next day's data point (2025-10-30).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1) "Prix Spot": insert the "01-nov" column before "01-oct." (col DF) ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# -4161 == xlShiftToRight: shift the existing DF:EJ columns right by one.
$wsPrix.Range("DF1:DF25").Insert(-4161)

# Header row keeps the same bold/centered/bordered look as its neighbours
# because Insert() carries the formatting of the following column along.
$wsPrix.Range("DF1").Value = "01-nov"

# No trading data yet for 01-nov, so every hour row gets a "-" placeholder,
# matching how other not-yet-available days are represented in this sheet.
for ($r = 2; $r -le 25; $r++) {
  $wsPrix.Cells.Item($r, 110).Value = "-"
}

# --- 2) "Gaz": append 2025-10-30 ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A138").NumberFormat = "@"
$wsGaz.Range("A138").Value = "2025-10-30"
$wsGaz.Range("A138").Style = "Normal"
$wsGaz.Range("B138").Value = 29.8

# --- 3) "CO2": append 2025-10-30 ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A138").NumberFormat = "@"
$wsCo2.Range("A138").Value = "2025-10-30"
$wsCo2.Range("A138").Style = "Normal"
$wsCo2.Range("B138").Value = 78.36
